# Update "想去人数" (F column) figures on the 展览 and 全部类型 sheets
# to reflect the latest generated output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of row number -> new F-column value
$updates = @{
    3  = 7094
    4  = 4832
    5  = 71
    10 = 69
    11 = 68
    12 = 189
    13 = 619
    14 = 136
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
